# Generate Report for Handback
#
# - Status changes from "Ready for handoff" to "Handed back: in sync with
#   en-US" on the Overview sheet and on each language sheet (zh-cn, de-de).
# - Each language sheet gets two new columns populated for the handback
#   report: F "Latest Target File" and G "Latest Handback File", both
#   rendered as hyperlinks (matching the look of the existing link columns).
# - The "Latest Handback DateTime" column (H) is stamped with the real
#   handback timestamp instead of the zero-date placeholder.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

function Set-HandbackRow {
    param($sheet, $row, $status, $mdUrl, $mdDisplay, $xlfUrl, $xlfDisplay, $handbackTime)

    $sheet.Range("C$row").Value = $status

    $sheet.Hyperlinks.Add($sheet.Range("F$row"), $mdUrl, "", "", $mdDisplay)
    $sheet.Hyperlinks.Add($sheet.Range("G$row"), $xlfUrl, "", "", $xlfDisplay)

    $sheet.Range("H$row").Value = $handbackTime
}

$zhcn = $wb.Worksheets.Item("zh-cn")

$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/92df28708ba510df4e54f68121a3f516bc2feef0/e2e/6ed8456b-c2ad-4330-9572-caf6d7fa80c0.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e854cb0ab8fbd4b68d4c74e0d076ea3e8f9e68e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6ed8456b-c2ad-4330-9572-caf6d7fa80c0.41d923ced29b8f299034b77d0df713481a4d1485.zh-cn.xlf"
$zhMdDisplay = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.md"
$zhXlfDisplay = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.41d923ced29b8f299034b77d0df713481a4d1485.zh-cn.xlf"

Set-HandbackRow $zhcn "2" $statusText $zhMdUrl $zhMdDisplay $zhXlfUrl $zhXlfDisplay "2016-03-18 07:29:19"
Set-HandbackRow $zhcn "3" $statusText $zhMdUrl $zhMdDisplay $zhXlfUrl $zhXlfDisplay "2016-03-18 07:29:19"

$dede = $wb.Worksheets.Item("de-de")

$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/92df28708ba510df4e54f68121a3f516bc2feef0/e2e/6ed8456b-c2ad-4330-9572-caf6d7fa80c0.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4e6478bb3d57a0f7b33fd616e5383d8a6f51147/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6ed8456b-c2ad-4330-9572-caf6d7fa80c0.41d923ced29b8f299034b77d0df713481a4d1485.de-de.xlf"
$deMdDisplay = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.md"
$deXlfDisplay = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.41d923ced29b8f299034b77d0df713481a4d1485.de-de.xlf"

Set-HandbackRow $dede "2" $statusText $deMdUrl $deMdDisplay $deXlfUrl $deXlfDisplay "2016-03-18 07:29:24"
Set-HandbackRow $dede "3" $statusText $deMdUrl $deMdDisplay $deXlfUrl $deXlfDisplay "2016-03-18 07:29:24"
